$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "9840022277"
Write-Host "done"
